# Update the 2024-06-12 cryptos snapshot: refresh Price/Volume(1h) figures
# and restore the Monero / FirstDigitalUSD row ordering (rows 40-41 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.381.34"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.502.98"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.87"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.27"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.502.23"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.05"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.093.59"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.507.95"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.372.01"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.14"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.20"
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.32"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.641.17"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +8.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.32"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -7.34%  "
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.09"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.493.04"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.03"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "177.71"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0870"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.878"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.37"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.46"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.25"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.987"
$ws.Range("E51").Value = "  -2.69%  "
